$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.741.37"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "2.583.85"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "519.59"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "139.17"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "2.594.03"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "6.53"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "3.035.87"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "58.734.96"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "20.38"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "2.562.81"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "338.27"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "10.13"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "6.49"
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "66.23"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").Value = "0.403"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "0.0₃0717"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("D31").Value = "5.92"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "18.75"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").Value = "148.54"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "36.49"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "274.74"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").Value = "10.75"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "0.0947"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "18.48"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("D49").Value = "1.976.82"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "4.49"
$ws.Range("E51").Value = "  -1.01%  "
